$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row number -> new Fitness (column C) value, per the commit
# "correction in sa algorithm and 746 logs". Only rows 2-149 are affected.
$newValues = @{
    2   = 8137
    3   = 8137
    4   = 8137
    5   = 8137
    6   = 8137
    7   = 8137
    8   = 8137
    9   = 8137
    10  = 8137
    11  = 8137
    12  = 8137
    13  = 8137
    14  = 8137
    15  = 8137
    16  = 8137
    17  = 8137
    18  = 8137
    19  = 8137
    20  = 8137
    21  = 8137
    22  = 8137
    23  = 8137
    24  = 8137
    25  = 8137
    26  = 8137
    27  = 8137
    28  = 8137
    29  = 8137
    30  = 8137
    31  = 8137
    32  = 8137
    33  = 8137
    34  = 8137
    35  = 8137
    36  = 8137
    37  = 8137
    38  = 8137
    39  = 8137
    40  = 8137
    41  = 8137
    42  = 8137
    43  = 8137
    44  = 8137
    45  = 8137
    46  = 8137
    47  = 8137
    48  = 8137
    49  = 8137
    50  = 8137
    51  = 8137
    52  = 8137
    53  = 8137
    54  = 8137
    55  = 8137
    56  = 8137
    57  = 8137
    58  = 8137
    59  = 8137
    60  = 8137
    61  = 8137
    62  = 8137
    63  = 8137
    64  = 8137
    65  = 8137
    66  = 8137
    67  = 8137
    68  = 8137
    69  = 8137
    70  = 8026
    71  = 8026
    72  = 7723
    73  = 7723
    74  = 7723
    75  = 7723
    76  = 7723
    77  = 7723
    78  = 7723
    79  = 7723
    80  = 7723
    81  = 7723
    82  = 7723
    83  = 7723
    84  = 7723
    85  = 7723
    86  = 7723
    87  = 7723
    88  = 7723
    89  = 7618
    90  = 7618
    91  = 7618
    92  = 7618
    93  = 7618
    94  = 7618
    95  = 7618
    96  = 7618
    97  = 7618
    98  = 7618
    99  = 7618
    100 = 7618
    101 = 7618
    102 = 7534
    103 = 7534
    104 = 7534
    105 = 7534
    106 = 7534
    107 = 7534
    108 = 7534
    109 = 7534
    110 = 7534
    111 = 7534
    112 = 7534
    113 = 7534
    114 = 7534
    115 = 7534
    116 = 7534
    117 = 7534
    118 = 7534
    119 = 7534
    120 = 7534
    121 = 7534
    122 = 7320
    123 = 7320
    124 = 7320
    125 = 7320
    126 = 7320
    127 = 7320
    128 = 7320
    129 = 7320
    130 = 7320
    131 = 7320
    132 = 7320
    133 = 7320
    134 = 7320
    135 = 7320
    136 = 7318
    137 = 7318
    138 = 7318
    139 = 7318
    140 = 7318
    141 = 7318
    142 = 7318
    143 = 7318
    144 = 7318
    145 = 7318
    146 = 7318
    147 = 7318
    148 = 7318
    149 = 7318
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
